$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.274.00"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").Value = "2.828.98"
$ws.Range("E3").Value = "  +1.04%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "356.54"
$ws.Range("E5").Value = "  +2.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "112.41"
$ws.Range("E6").Value = "  -3.70%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.570"
$ws.Range("E7").Value = "  +3.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  +1.84%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.97"
$ws.Range("E10").Value = "  -5.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0864"
$ws.Range("E11").Value = "  +0.70%  "
$ws.Range("E12").Value = "  +0.92%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.96"
$ws.Range("E13").Value = "  -0.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.78"
$ws.Range("E14").Value = "  -0.87%  "
$ws.Range("D15").Value = "3.269.95"
$ws.Range("E15").Value = "  +0.90%  "
$ws.Range("D16").Value = "2.825.02"
$ws.Range("E16").Value = "  +1.64%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.933"
$ws.Range("E17").Value = "  +4.26%  "
$ws.Range("D18").Value = "52.102.80"
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.54"
$ws.Range("E19").Value = "  +5.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.21"
$ws.Range("E20").Value = "  -0.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.44"
$ws.Range("E21").Value = "  -0.45%  "
$ws.Range("D22").Value = "0.0₃0998"
$ws.Range("E22").Value = "  +1.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.88"
$ws.Range("E23").Value = "  +0.94%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "272.18"
$ws.Range("E24").Value = "  +0.75%  "
$ws.Range("E25").Value = "  +2.62%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "27.04"
$ws.Range("E26").Value = "  +1.09%  "
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("E28").Value = "  +0.85%  "
$ws.Range("E29").Value = "  +0.74%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.146"
$ws.Range("E30").Value = "  +3.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0488"
$ws.Range("E31").Value = "  +17.92%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "52.46"
$ws.Range("E32").Value = "  +4.37%  "
$ws.Range("E33").Value = "  -0.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.95"
$ws.Range("E34").Value = "  +4.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.62"
$ws.Range("E35").Value = "  +12.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0857"
$ws.Range("E36").Value = "  +3.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.29"
$ws.Range("E38").Value = "  +2.06%  "
$ws.Range("E39").Value = "  -3.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.37"
$ws.Range("E40").Value = "  -2.88%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.117"
$ws.Range("E41").Value = "  +1.64%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "127.48"
$ws.Range("E42").Value = "  -0.76%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "23.35"
$ws.Range("E43").Value = "  -0.77%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.53"
$ws.Range("E44").Value = "  -6.93%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.27"
$ws.Range("E45").Value = "  -1.70%  "
$ws.Range("E46").Value = "  +0.49%  "
$ws.Range("D47").Value = "2.091.19"
$ws.Range("E47").Value = "  +0.81%  "
$ws.Range("E48").Value = "  -4.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.93"
$ws.Range("E49").Value = "  +7.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.982"
$ws.Range("E50").Value = "  -0.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.19"
$ws.Range("E51").Value = "  +2.48%  "
